$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "OD" header row; everything below shifts up one row.
$ws.Rows.Item(1).Delete()

# Match the new selection recorded after the edit.
$ws.Range("D9").Select()
